# atualizei dados bibi e add
# Insert a new "venda atipica" record for 2025-06-30 (id_venda 374455,
# "CAIXA DE SOM SEM FIO 5W METAL") ahead of the existing 2025-06-30 rows,
# and refresh the computed stats (estoque/media/desvio) for the rows that
# follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: push current rows 3-8 down to 4-9.
$ws.Rows.Item(3).Insert()

function Set-TextCell($cell, $text) {
    # Force text storage for values that look numeric (ids, dates) so they
    # keep being shared-string cells instead of being coerced to numbers,
    # then drop the temporary "@" number format so no stray style sticks.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# New row 3: 2025-06-30 / BEMOL S/A / id_venda 374455 / CAIXA DE SOM SEM FIO 5W METAL
Set-TextCell "A3" "2025-06-30"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "BEMOL S/A"
Set-TextCell "D3" "374455"
$ws.Range("E3").Value = 47791
Set-TextCell "F3" "CAIXA DE SOM SEM FIO 5W METAL"
$ws.Range("G3").Value = -153
$ws.Range("H3").Value = 1.06
$ws.Range("I3").Value = 0.28

# Refresh recalculated columns (E, G, H, I) on the rows that shifted down.
$ws.Range("E4").Value = 13079
$ws.Range("G4").Value = -462
$ws.Range("H4").Value = 1.08
$ws.Range("I4").Value = 0.28

$ws.Range("E5").Value = 49904
$ws.Range("G5").Value = -17
$ws.Range("H5").Value = 1.06
$ws.Range("I5").Value = 0.25

$ws.Range("E6").Value = 14186
$ws.Range("G6").Value = -6
$ws.Range("H6").Value = 1.03
$ws.Range("I6").Value = 0.18

$ws.Range("E7").Value = 10114
$ws.Range("G7").Value = -85
$ws.Range("H7").Value = 1.05
$ws.Range("I7").Value = 0.22

$ws.Range("E8").Value = 13018
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1.08
$ws.Range("I8").Value = 0.29

$ws.Range("E9").Value = 13546
$ws.Range("G9").Value = -322
$ws.Range("H9").Value = 1.1
$ws.Range("I9").Value = 0.34
